$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.792.81'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '3.264.87'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.46'
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.65'
$ws.Range("E6").Value = '  +0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("E9").Value = '  -2.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.59'
$ws.Range("E10").Value = '  -1.21%  '
$ws.Range("E11").Value = '  -4.15%  '
$ws.Range("D12").Value = '3.834.11'
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.34'
$ws.Range("E14").Value = '  -4.48%  '
$ws.Range("D15").Value = '67.794.89'
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("E16").Value = '  -2.37%  '
$ws.Range("D17").Value = '3.238.72'
$ws.Range("E17").Value = '  +1.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.68'
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.38'
$ws.Range("E19").Value = '  -1.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '401.98'
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.53'
$ws.Range("E21").Value = '  -2.31%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.05'
$ws.Range("E23").Value = '  -1.37%  '
$ws.Range("E24").Value = '  -1.84%  '
$ws.Range("E25").Value = '  -2.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.187'
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.45'
$ws.Range("E27").Value = '  -2.16%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.63'
$ws.Range("E30").Value = '  -1.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.45'
$ws.Range("E31").Value = '  -4.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.88'
$ws.Range("E32").Value = '  -3.98%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  -3.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '164.54'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.45'
$ws.Range("E36").Value = '  -3.90%  '
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.10'
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.801'
$ws.Range("E39").Value = '  -3.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.47'
$ws.Range("E40").Value = '  -3.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.34'
$ws.Range("E41").Value = '  -3.09%  '
$ws.Range("D42").Value = '2.680.52'
$ws.Range("E42").Value = '  +2.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.78'
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0676'
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("E45").Value = '  -3.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '334.35'
$ws.Range("E46").Value = '  -3.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.54'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("E48").Value = '  -3.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.29'
$ws.Range("E49").Value = '  -0.74%  '
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.965'
$ws.Range("E51").Value = '  -1.94%  '
